# The deck's Design was switched from the "Integral" theme to the
# built-in "Office Theme" (the two <a:theme> parts trade places: the
# palette that used to be on the slide master's theme becomes the
# classic Office blue palette, and vice-versa). Re-create that by
# pushing the Office Theme's 12 theme colors onto the presentation's
# live theme color scheme (the one actually used by the slide master /
# slides), the same way PowerPoint itself rewrites the colours when you
# click a different swatch in the Design gallery.

$p = $ppt.ActivePresentation

function RGBVal([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme colour scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$tcs = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 0; $i -lt $officeColors.Count; $i++) {
    $tcs.Colors($i + 1).RGB = RGBVal($officeColors[$i])
}
